$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.178.06"
$ws.Range("E2").Value = "  +2.13%  "
$ws.Range("D3").Value = "1.980.94"
$ws.Range("E3").Value = "  +4.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9967"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8008"
$ws.Range("E5").Value = "  +69.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "253.72"
$ws.Range("E6").Value = "  +3.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9969"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E8").Value = "  +18.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.69"
$ws.Range("E9").Value = "  +15.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06938"
$ws.Range("E10").Value = "  +6.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8533"
$ws.Range("E11").Value = "  +17.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08149"
$ws.Range("E12").Value = "  +4.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "102.41"
$ws.Range("E13").Value = "  +6.77%  "
$ws.Range("D14").Value = "1.975.33"
$ws.Range("E14").Value = "  +4.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.501"
$ws.Range("E15").Value = "  +5.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "279.63"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").Value = "31.146.33"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.95"
$ws.Range("E18").Value = "  +6.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007878"
$ws.Range("E19").Value = "  +5.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.690"
$ws.Range("E20").Value = "  +7.97%  "
$ws.Range("D21").Value = "2.226.24"
$ws.Range("E21").Value = "  +4.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9971"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9955"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.809"
$ws.Range("E24").Value = "  +7.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1590"
$ws.Range("E25").Value = "  +64.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.691"
$ws.Range("E26").Value = "  +6.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.30"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.63"
$ws.Range("E28").Value = "  +3.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.201"
$ws.Range("E29").Value = "  +16.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.564"
$ws.Range("E30").Value = "  +6.33%  "
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.560"
$ws.Range("E32").Value = "  +6.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.353"
$ws.Range("E33").Value = "  +4.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05144"
$ws.Range("E34").Value = "  +5.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.223"
$ws.Range("E35").Value = "  +8.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7405"
$ws.Range("E36").Value = "  +6.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.803"
$ws.Range("E37").Value = "  +3.22%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01989"
$ws.Range("E38").Value = "  +5.36%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.909"
$ws.Range("E39").Value = "  +3.08%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.600"
$ws.Range("E40").Value = "  +6.33%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "78.42"
$ws.Range("E41").Value = "  +4.73%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4695"
$ws.Range("E42").Value = "  +9.90%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.078"
$ws.Range("E43").Value = "  +5.46%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.26"
$ws.Range("E44").Value = "  +3.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8497"
$ws.Range("E45").Value = "  +2.70%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9968"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.960"
$ws.Range("E47").Value = "  +3.21%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.538"
$ws.Range("E48").Value = "  +8.33%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4267"
$ws.Range("E49").Value = "  +8.48%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.33"
$ws.Range("E50").Value = "  +3.27%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "934.01"
$ws.Range("E51").Value = "  +2.96%  "
